# "CP terminados, actualizado el script de precondiciones"
#
# Caso de prueba 9 "Registrar Usuario" -> Caso de prueba 11 "Crear Rol":
#  - DatosGenerales: ID/Nombre/Descripcion updated
#  - Precondiciones: precondition text updated + a second precondition row added
#  - Pasos: steps 2-4 filled in with the "Crear Rol" flow, trailing empty step
#    rows (5 & 6) lose their numbering
#  - DatosGenerales active tab moves back to the first sheet

$wb = $excel.ActiveWorkbook

$wsDatos   = $wb.Worksheets.Item("DatosGenerales")
$wsPrecond = $wb.Worksheets.Item("Precondiciones")
$wsPasos   = $wb.Worksheets.Item("Pasos")

# ---------------------------------------------------------------------------
# DatosGenerales: new ID / Nombre / Descripcion for the test case
# ---------------------------------------------------------------------------
$wsDatos.Range("B1").Value = "11"
$wsDatos.Range("B2").Value = "Crear Rol"
$wsDatos.Range("B3").Value = "Crear un nuevo rol en el sistema"

# ---------------------------------------------------------------------------
# Precondiciones: rewrite the existing precondition and add a new one
# ---------------------------------------------------------------------------
$wsPrecond.Range("B2").Value = 'El rol "CPA_Rol" no existe en la base de datos'
$wsPrecond.Rows.Item(2).AutoFit()

$wsPrecond.Range("A3").Value = 2
$wsPrecond.Range("B3").Value = "El usuario logueado es un usuario administrador"

# ---------------------------------------------------------------------------
# Pasos: fill in steps 2-4 of the "Crear Rol" flow
# ---------------------------------------------------------------------------
$wsPasos.Range("B3").Value = "Seleccionar la opcion Administracon de usuarios > Nuevo Rol"
$wsPasos.Range("C3").Value = "Se muestra la pagina AdministracionUsuarios.aspx?action=NuevoRol"
$wsPasos.Rows.Item(3).RowHeight = 26.25

$wsPasos.Range("B4").Value = 'Ingreso "CPA_Rol" en el campo nombre y "CPA_Descripcion" en el campo descripcion'
$wsPasos.Rows.Item(4).RowHeight = 26.25

$wsPasos.Range("B5").Value = "Presiono el boton Guardar"
$wsPasos.Range("C5").Value = "Se muestra un mensaje de exito."

# Steps 5 & 6 are no longer used - clear their step numbers but keep formatting
$wsPasos.Range("A6").ClearContents()
$wsPasos.Range("A7").ClearContents()

# ---------------------------------------------------------------------------
# Selections / active sheet: DatosGenerales becomes the active tab again
# ---------------------------------------------------------------------------
$wsPrecond.Range("A4").Select() | Out-Null
$wsPasos.Range("C7").Select() | Out-Null

$wsDatos.Activate()
$wsDatos.Range("C14").Select() | Out-Null
